# no-op test
